$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update input parameters / recalculated figures for the investment
# --- calculation ("class format" update referenced in the commit message).
# Grundinvestering (initial investment) halved
$ws.Range("B2").Value = -1000000

# Avskrivningar (depreciation), 3% of the new investment, years 1-10
$ws.Range("C3").Value = 30000
$ws.Range("D3").Value = 30000
$ws.Range("E3").Value = 30000
$ws.Range("F3").Value = 30000
$ws.Range("G3").Value = 30000
$ws.Range("H3").Value = 30000
$ws.Range("I3").Value = 30000
$ws.Range("J3").Value = 30000
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 30000

# Inbetalningar (income), years 1-10
$ws.Range("C4").Value = 489999.9999999999
$ws.Range("D4").Value = 489999.9999999999
$ws.Range("E4").Value = 489999.9999999999
$ws.Range("F4").Value = 489999.9999999999
$ws.Range("G4").Value = 489999.9999999999
$ws.Range("H4").Value = 489999.9999999999
$ws.Range("I4").Value = 489999.9999999999
$ws.Range("J4").Value = 489999.9999999999
$ws.Range("K4").Value = 489999.9999999999
$ws.Range("L4").Value = 489999.9999999999

# Utbetalningar (expenses), years 1-10
$ws.Range("C5").Value = -140000
$ws.Range("D5").Value = -140000
$ws.Range("E5").Value = -140000
$ws.Range("F5").Value = -140000
$ws.Range("G5").Value = -140000
$ws.Range("H5").Value = -140000
$ws.Range("I5").Value = -140000
$ws.Range("J5").Value = -140000
$ws.Range("K5").Value = -140000
$ws.Range("L5").Value = -140000

# Rorelsebindandekapital (working capital), tripled
$ws.Range("B7").Value = -300000
$ws.Range("L7").Value = 300000

# Arligt netto (annual net)
$ws.Range("B8").Value = -1440000
$ws.Range("C8").Value = 379999.9999999999
$ws.Range("D8").Value = 379999.9999999999
$ws.Range("E8").Value = 379999.9999999999
$ws.Range("F8").Value = 379999.9999999999
$ws.Range("G8").Value = 379999.9999999999
$ws.Range("H8").Value = 379999.9999999999
$ws.Range("I8").Value = 379999.9999999999
$ws.Range("J8").Value = 379999.9999999999
$ws.Range("K8").Value = 379999.9999999999
$ws.Range("L8").Value = 820000

# Nuvarde (present value)
$ws.Range("B9").Value = -1440000
$ws.Range("C9").Value = 343891.4027149321
$ws.Range("D9").Value = 311213.9391085358
$ws.Range("E9").Value = 281641.5738538786
$ws.Range("F9").Value = 254879.2523564512
$ws.Range("G9").Value = 230659.9568836663
$ws.Range("H9").Value = 208742.0424286572
$ws.Range("I9").Value = 188906.8257272916
$ws.Range("J9").Value = 170956.4033731146
$ws.Range("K9").Value = 154711.6772607372
$ws.Range("L9").Value = 302128.0670488319

# Ackumelerat nuvarde (accumulated present value)
$ws.Range("B10").Value = -1440000
$ws.Range("C10").Value = -1096108.597285068
$ws.Range("D10").Value = -784894.658176532
$ws.Range("E10").Value = -503253.0843226534
$ws.Range("F10").Value = -248373.8319662022
$ws.Range("G10").Value = -17713.87508253596
$ws.Range("H10").Value = 191028.1673461213
$ws.Range("I10").Value = 379934.9930734129
$ws.Range("J10").Value = 550891.3964465274
$ws.Range("K10").Value = 705603.0737072646
$ws.Range("L10").Value = 1007731.140756096

# Nettonuvarde (net present value)
$ws.Range("B11").Value = 1007731.140756096

# The accumulated present value in F10/G10 turned negative under the new
# figures, so re-apply the "negative" (red) formatting used elsewhere in the
# column instead of the "positive" (green) one it previously had. Copy the
# cell format (fill/number format) from B5, which already carries the
# negative-value style, onto F10 and G10.
$ws.Range("B5").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
